$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 1.199219
$ws.Range("N2").Value = 3.597657
$ws.Range("O2").Value = 0.1760050710731031
$ws.Range("P2").Value = 0.1760050710731032
$ws.Range("Q2").Value = 0.07804197460266667
$ws.Range("R2").Value = 0.7023777714239999
$ws.Range("S2").Value = 0.1760050710731031
$ws.Range("T2").Value = 0.1760050710731032

# Row 3
$ws.Range("O3").Value = 0.5218760230842041
$ws.Range("P3").Value = 0.5218760230842042
$ws.Range("S3").Value = 0.5218760230842041
$ws.Range("T3").Value = 0.5218760230842042

# Row 4
$ws.Range("M4").Value = 2.058501666666667
$ws.Range("N4").Value = 6.175505
$ws.Range("O4").Value = 0.3021189058426926
$ws.Range("P4").Value = 0.3021189058426926
$ws.Range("S4").Value = 0.3021189058426926
$ws.Range("T4").Value = 0.3021189058426926
